# edit.ps1
# Applies the data-correction edit described by the commit "Add files via upload":
#   - Updated the CWON/Fragile-States indicator figures in rows 97 ("Africa, Fragile
#     States") and 98 ("ROW, Fragile States") of the Tab29 worksheet.
#   - Fixed mis-encoded accented characters (mojibake) in the "Regional Economic
#     Communities" footnote so that PALOP's definition reads correctly in
#     Portuguese/Spanish ("Países Africanos de Língua Oficial Portuguesa" /
#     "Mercado Común del Sur").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tab29")

# --- Row 97: "Africa, Fragile States" updated indicator values ---
$ws.Range("C97").Value = 12868.541997093
$ws.Range("D97").Value = 3272.5021003460001
$ws.Range("E97").Value = 3290.658625392
$ws.Range("F97").Value = 8141.6481835340001
$ws.Range("G97").Value = 2019.0304741949999
$ws.Range("H97").Value = 395.700717798911
$ws.Range("I97").Value = 304.22211245789998
$ws.Range("J97").Value = 5.0475930165899996
$ws.Range("K97").Value = 10.457351998115
$ws.Range("L97").Value = 233.34297219338299
$ws.Range("M97").Value = 1070.2603219179
$ws.Range("N97").Value = 717.83695006790003
$ws.Range("O97").Value = 352.4235396222
$ws.Range("P97").Value = 1271.6267689962399
$ws.Range("Q97").Value = 1201.9078213620401
$ws.Range("R97").Value = 1057.4842536707399
$ws.Range("S97").Value = 132.1887516834
$ws.Range("U97").Value = 69.719462086830006
$ws.Range("V97").Value = -245.56838999999999

# --- Row 98: "ROW, Fragile States" updated indicator values ---
$ws.Range("C98").Value = 22835.330017565
$ws.Range("D98").Value = 5399.5377050799998
$ws.Range("E98").Value = 4855.6874826869998
$ws.Range("F98").Value = 13047.396396308
$ws.Range("G98").Value = 1993.3771594150001
$ws.Range("H98").Value = 119.17308086803
$ws.Range("I98").Value = 276.51773982840001
$ws.Range("J98").Value = 19.60991159316
$ws.Range("K98").Value = 13.3793097011
$ws.Range("L98").Value = 218.68722159559999
$ws.Range("M98").Value = 1346.010055345
$ws.Range("N98").Value = 880.31737984699998
$ws.Range("O98").Value = 465.69361021934998
$ws.Range("P98").Value = 2917.3534394447802
$ws.Range("Q98").Value = 2835.3682857867798
$ws.Range("R98").Value = 2344.1877478162301
$ws.Range("S98").Value = 485.24530746568001
$ws.Range("T98").Value = 5.9257335124840003
$ws.Range("U98").Value = 81.992282101800001
$ws.Range("V98").Value = 17.601220000000001

# --- Footnote text fix (A103): correct mojibake accented characters ---
$ws.Range("A103").Value = 'Regional Economic Communities:CEN-SAD = "Community of Sahel-Saharan States";COMESA = "Common Market for Eastern and Southern Africa";EAC = "East African Community";ECCAS = "Economic Community of Central African States";ECOWAS = "Economic Community of West African States";IGAD = "Intergovernmental Authority on Development";SADC = "Southern African Development Community";UMA = "Arab Maghreb Union";PALOP = "Países Africanos de Língua Oficial Portuguesa";ASEAN = "Association of Southeast Asian Nations";MERCOSUR = "Mercado Común del Sur".EU27 = "European Union (27 members)".OECD = "Organisation for Economic Co-operation and Development".'
